$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 1017.25
$ws.Cells.Item(32, 9).Value = 499
$ws.Cells.Item(32, 10).Value = 1190
$ws.Cells.Item(32, 11).Value = 499
$ws.Cells.Item(32, 12).Value = 1190
$ws.Cells.Item(32, 13).Value = -173
$ws.Cells.Item(32, 14).Value = -1842

$ws.Cells.Item(74, 8).Value = 3704
$ws.Cells.Item(74, 9).Value = 3350.75
$ws.Cells.Item(74, 10).Value = 4410.5
$ws.Cells.Item(74, 11).Value = 3350.75
$ws.Cells.Item(74, 12).Value = 4410.5
$ws.Cells.Item(74, 13).Value = -2414.75
$ws.Cells.Item(74, 14).Value = -6282.5

$ws.Cells.Item(76, 8).Value = 2472848.5
$ws.Cells.Item(76, 9).Value = 2852589.2
$ws.Cells.Item(76, 10).Value = 4533.5
$ws.Cells.Item(76, 11).Value = 2852589.2
$ws.Cells.Item(76, 12).Value = 4533.5
$ws.Cells.Item(76, 13).Value = -2852274.2
$ws.Cells.Item(76, 14).Value = -5163.5

$ws.Cells.Item(77, 8).Value = 3704
$ws.Cells.Item(77, 9).Value = 3350.75
$ws.Cells.Item(77, 10).Value = 4410.5
$ws.Cells.Item(77, 11).Value = 16753.75
$ws.Cells.Item(77, 12).Value = 22052.5
$ws.Cells.Item(77, 13).Value = -12073.75
$ws.Cells.Item(77, 14).Value = -31412.5

$ws.Cells.Item(79, 8).Value = 2472848.5
$ws.Cells.Item(79, 9).Value = 2852589.2
$ws.Cells.Item(79, 10).Value = 4533.5
$ws.Cells.Item(79, 11).Value = 2852589.2
$ws.Cells.Item(79, 12).Value = 4533.5
$ws.Cells.Item(79, 13).Value = -2851497.2
$ws.Cells.Item(79, 14).Value = -6717.5

$ws.Cells.Item(112, 8).Value = 5929.909
$ws.Cells.Item(112, 10).Value = 6874.5356
$ws.Cells.Item(112, 12).Value = 20623.6068
$ws.Cells.Item(112, 14).Value = -22839.6068

$ws.Cells.Item(137, 8).Value = 1860.8235
$ws.Cells.Item(137, 9).Value = 1331.875
$ws.Cells.Item(137, 11).Value = 3995.625
$ws.Cells.Item(137, 13).Value = -1445.625

$ws.Cells.Item(138, 8).Value = 4985.0684
$ws.Cells.Item(138, 10).Value = 5342.067
$ws.Cells.Item(138, 12).Value = 16026.201
$ws.Cells.Item(138, 14).Value = -26306.201

$ws.Cells.Item(140, 8).Value = 97966.664
$ws.Cells.Item(140, 10).Value = 97966.664
$ws.Cells.Item(140, 12).Value = 97966.664
$ws.Cells.Item(140, 14).Value = -108326.664

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 608722.9399999999
$ws.Cells.Item(32, 9).Value = 10094.0625
$ws.Cells.Item(32, 10).Value = 1758090.4
$ws.Cells.Item(32, 11).Value = 10094.0625
$ws.Cells.Item(32, 12).Value = 1758090.4
$ws.Cells.Item(32, 13).Value = -9807.0625
$ws.Cells.Item(32, 14).Value = -1758664.4

$ws.Cells.Item(44, 8).Value = 6857.143
$ws.Cells.Item(44, 10).Value = 6857.143
$ws.Cells.Item(44, 12).Value = 6857.143
$ws.Cells.Item(44, 14).Value = -7833.143

$ws.Cells.Item(61, 8).Value = 3462.4666
$ws.Cells.Item(61, 9).Value = 3608.3845
$ws.Cells.Item(61, 11).Value = 3608.3845
$ws.Cells.Item(61, 13).Value = -3396.3845

$ws.Cells.Item(63, 8).Value = 8000
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 10).Value = 8000
$ws.Cells.Item(63, 11).Value = 0
$ws.Cells.Item(63, 12).Value = 8000
$ws.Cells.Item(63, 14).Value = -9372
$ws.Cells.Item(63, 13).ClearContents()

$ws.Cells.Item(66, 8).Value = 8000
$ws.Cells.Item(66, 9).Value = 0
$ws.Cells.Item(66, 10).Value = 8000
$ws.Cells.Item(66, 11).Value = 0
$ws.Cells.Item(66, 12).Value = 40000
$ws.Cells.Item(66, 14).Value = -46864
$ws.Cells.Item(66, 13).ClearContents()

$ws.Cells.Item(74, 8).Value = 1057.8
$ws.Cells.Item(74, 9).Value = 1219
$ws.Cells.Item(74, 10).Value = 413
$ws.Cells.Item(74, 11).Value = 1219
$ws.Cells.Item(74, 12).Value = 413
$ws.Cells.Item(74, 13).Value = -345
$ws.Cells.Item(74, 14).Value = -2161

$ws.Cells.Item(77, 8).Value = 1057.8
$ws.Cells.Item(77, 9).Value = 1219
$ws.Cells.Item(77, 10).Value = 413
$ws.Cells.Item(77, 11).Value = 6095
$ws.Cells.Item(77, 12).Value = 2065
$ws.Cells.Item(77, 13).Value = -1727
$ws.Cells.Item(77, 14).Value = -10801

$ws.Cells.Item(132, 8).Value = 2765.5527
$ws.Cells.Item(132, 9).Value = 2358.1428
$ws.Cells.Item(132, 10).Value = 3906.3
$ws.Cells.Item(132, 11).Value = 7074.428400000001
$ws.Cells.Item(132, 12).Value = 11718.9
$ws.Cells.Item(132, 13).Value = -4544.428400000001
$ws.Cells.Item(132, 14).Value = -16778.9

$ws.Cells.Item(136, 8).Value = 3462.4666
$ws.Cells.Item(136, 9).Value = 3608.3845
$ws.Cells.Item(136, 11).Value = 10825.1535
$ws.Cells.Item(136, 13).Value = -8275.1535

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 5133.3335
$ws.Cells.Item(105, 9).Value = 5133.3335
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 11).Value = 5133.3335
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 13).Value = -3386.3335
$ws.Cells.Item(105, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 9755.777
$ws.Cells.Item(31, 9).Value = 3255.7568
$ws.Cells.Item(31, 11).Value = 3255.7568
$ws.Cells.Item(31, 13).Value = -2960.7568

$ws.Cells.Item(34, 8).Value = 9755.777
$ws.Cells.Item(34, 9).Value = 3255.7568
$ws.Cells.Item(34, 11).Value = 3255.7568
$ws.Cells.Item(34, 13).Value = -3053.7568

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(58, 8).Value = 984.8333
$ws.Cells.Item(58, 9).Value = 436.33334
$ws.Cells.Item(58, 10).Value = 1533.3334
$ws.Cells.Item(58, 11).Value = 1309.00002
$ws.Cells.Item(58, 12).Value = 4600.0002
$ws.Cells.Item(58, 13).Value = -1181.00002
$ws.Cells.Item(58, 14).Value = -4856.0002

$ws.Cells.Item(107, 8).Value = 485.58066
$ws.Cells.Item(107, 9).Value = 255.36842
$ws.Cells.Item(107, 10).Value = 850.0833
$ws.Cells.Item(107, 11).Value = 766.1052599999999
$ws.Cells.Item(107, 12).Value = 2550.2499
$ws.Cells.Item(107, 13).Value = 1153.89474
$ws.Cells.Item(107, 14).Value = -6390.2499

$ws.Cells.Item(113, 8).Value = 884.02
$ws.Cells.Item(113, 9).Value = 473
$ws.Cells.Item(113, 10).Value = 929.6889
$ws.Cells.Item(113, 11).Value = 1419
$ws.Cells.Item(113, 12).Value = 2789.0667
$ws.Cells.Item(113, 13).Value = 751
$ws.Cells.Item(113, 14).Value = -7129.066699999999

$ws.Cells.Item(121, 8).Value = 107570.266
$ws.Cells.Item(121, 9).Value = 3283.6
$ws.Cells.Item(121, 10).Value = 223444.33
$ws.Cells.Item(121, 11).Value = 9850.799999999999
$ws.Cells.Item(121, 12).Value = 670332.99
$ws.Cells.Item(121, 13).Value = -8540.799999999999
$ws.Cells.Item(121, 14).Value = -672952.99

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 25577454
$ws.Cells.Item(70, 9).Value = 40186644
$ws.Cells.Item(70, 10).Value = 11375
$ws.Cells.Item(70, 11).Value = 40186644
$ws.Cells.Item(70, 12).Value = 11375
$ws.Cells.Item(70, 13).Value = -40186374
$ws.Cells.Item(70, 14).Value = -11915

$ws.Cells.Item(73, 8).Value = 25577454
$ws.Cells.Item(73, 9).Value = 40186644
$ws.Cells.Item(73, 10).Value = 11375
$ws.Cells.Item(73, 11).Value = 40186644
$ws.Cells.Item(73, 12).Value = 11375
$ws.Cells.Item(73, 13).Value = -40185708
$ws.Cells.Item(73, 14).Value = -13247

$ws.Cells.Item(80, 8).Value = 4330
$ws.Cells.Item(80, 9).Value = 4995
$ws.Cells.Item(80, 10).Value = 3000
$ws.Cells.Item(80, 11).Value = 4995
$ws.Cells.Item(80, 12).Value = 3000
$ws.Cells.Item(80, 13).Value = -3997
$ws.Cells.Item(80, 14).Value = -4996

$ws.Cells.Item(83, 8).Value = 4330
$ws.Cells.Item(83, 9).Value = 4995
$ws.Cells.Item(83, 10).Value = 3000
$ws.Cells.Item(83, 11).Value = 24975
$ws.Cells.Item(83, 12).Value = 15000
$ws.Cells.Item(83, 13).Value = -19983
$ws.Cells.Item(83, 14).Value = -24984

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 1043.1471
$ws.Cells.Item(136, 10).Value = 2857.8572
$ws.Cells.Item(136, 12).Value = 8573.571599999999
$ws.Cells.Item(136, 14).Value = -13673.5716
